# Update the timetable sheet to work with "groups" instead of day/time slots.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rename the Lab slot headers, bold them ---
$ws.Range("C1").Value = "Lab 1"
$ws.Range("D1").Value = "Lab 2"
$ws.Range("E1").Value = "Lab 3"
$ws.Range("F1").Value = "Lab 4"
$ws.Range("G1").Value = "Lab 5"

$headerRange = $ws.Range("A1:G1")
$headerRange.Font.Bold = $true

# --- Body values (rows 2-6, columns C-G): replace lab rotation with group rotation ---
$groups = @(
    @("Group A", "Group E", "Group D", "Group C", "Group B"),
    @("Group B", "Group A", "Group E", "Group D", "Group C"),
    @("Group C", "Group B", "Group A", "Group E", "Group D"),
    @("Group D", "Group C", "Group B", "Group A", "Group E"),
    @("Group E", "Group D", "Group C", "Group B", "Group A")
)

for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $rowVals = $groups[$i]
    for ($j = 0; $j -lt 5; $j++) {
        $col = 3 + $j
        $ws.Cells.Item($row, $col).Value = $rowVals[$j]
    }
    # shrink the row height now that the text is shorter than the old day/time labels
    $ws.Rows.Item($row).RowHeight = 16.5
}

# match the number-formatted, right-aligned / wrapped style used by the rest of the row (B column)
$ws.Range("B2").Copy()
$ws.Range("C2:G6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column widths for the new group columns ---
$ws.Range("C1:G6").EntireColumn.AutoFit()

# --- Selection, matching the saved view state ---
$ws.Range("C2:G6").Select()

$wb.Save()
